$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template values shared by every row in this block (Mercado/Región/Tipo/Producto/Categoría/Kg per unit)
$mercadoId   = 1
$mercado     = "Agrícola del Norte S.A. de Arica"
$region      = "Arica y Parinacota"
$codreg      = 15
$tipo        = "Fruta"
$productoId  = 100104
$producto    = "Frutos de pepita"
$categoriaId = 100104002
$categoria   = "Manzana"
$kgUnidad    = 18
$fechaFormat = $ws.Cells.Item(49, 4).NumberFormat

# Row 49
$ws.Cells.Item(49, 4).Value = 44435
$ws.Cells.Item(49, 11).Value = "Granny Smith"
$ws.Cells.Item(49, 12).Value = "Calibre 90"
$ws.Cells.Item(49, 13).Value = 300
$ws.Cells.Item(49, 14).Value = 17000
$ws.Cells.Item(49, 15).Value = 18000
$ws.Cells.Item(49, 16).Value = 17500
$ws.Cells.Item(49, 17).Value = "$/caja 18 kilos embalada"
$ws.Cells.Item(49, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(49, 19).Value = 972

# Row 50
$ws.Cells.Item(50, 4).Value = 44435
$ws.Cells.Item(50, 11).Value = "Pink Lady"
$ws.Cells.Item(50, 12).Value = "Calibre 80"
$ws.Cells.Item(50, 13).Value = 250
$ws.Cells.Item(50, 14).Value = 17000
$ws.Cells.Item(50, 15).Value = 18000
$ws.Cells.Item(50, 16).Value = 17500
$ws.Cells.Item(50, 17).Value = "$/caja 18 kilos embalada"
$ws.Cells.Item(50, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(50, 19).Value = 972

# Row 51
$ws.Cells.Item(51, 4).Value = 44435
$ws.Cells.Item(51, 11).Value = "Royal Gala"
$ws.Cells.Item(51, 12).Value = "Calibre 90"
$ws.Cells.Item(51, 13).Value = 520
$ws.Cells.Item(51, 14).Value = 17000
$ws.Cells.Item(51, 15).Value = 18000
$ws.Cells.Item(51, 16).Value = 17500
$ws.Cells.Item(51, 17).Value = "$/caja 18 kilos embalada"
$ws.Cells.Item(51, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(51, 19).Value = 972

# Row 52
$ws.Cells.Item(52, 4).Value = 44435
$ws.Cells.Item(52, 11).Value = "Scarlett"
$ws.Cells.Item(52, 12).Value = "Calibre 80"
$ws.Cells.Item(52, 13).Value = 300
$ws.Cells.Item(52, 14).Value = 17000
$ws.Cells.Item(52, 15).Value = 18000
$ws.Cells.Item(52, 16).Value = 17500
$ws.Cells.Item(52, 17).Value = "$/caja 18 kilos embalada"
$ws.Cells.Item(52, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(52, 19).Value = 972

# Row 53
$ws.Cells.Item(53, 4).Value = 44351
$ws.Cells.Item(53, 11).Value = "Fuji royal"
$ws.Cells.Item(53, 12).Value = "Segunda"
$ws.Cells.Item(53, 13).Value = 300
$ws.Cells.Item(53, 14).Value = 16000
$ws.Cells.Item(53, 15).Value = 17000
$ws.Cells.Item(53, 16).Value = 16500
$ws.Cells.Item(53, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(53, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(53, 19).Value = 917

# Row 54
$ws.Cells.Item(54, 4).Value = 44351
$ws.Cells.Item(54, 11).Value = "Granny Smith"
$ws.Cells.Item(54, 12).Value = "Segunda"
$ws.Cells.Item(54, 13).Value = 300
$ws.Cells.Item(54, 14).Value = 16000
$ws.Cells.Item(54, 15).Value = 17000
$ws.Cells.Item(54, 16).Value = 16500
$ws.Cells.Item(54, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(54, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(54, 19).Value = 917

# Row 55
$ws.Cells.Item(55, 4).Value = 44351
$ws.Cells.Item(55, 11).Value = "Royal Gala"
$ws.Cells.Item(55, 12).Value = "Segunda"
$ws.Cells.Item(55, 13).Value = 250
$ws.Cells.Item(55, 14).Value = 16000
$ws.Cells.Item(55, 15).Value = 17000
$ws.Cells.Item(55, 16).Value = 16500
$ws.Cells.Item(55, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(55, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(55, 19).Value = 917

# Row 56
$ws.Cells.Item(56, 4).Value = 44278
$ws.Cells.Item(56, 11).Value = "Fuji royal"
$ws.Cells.Item(56, 12).Value = "Segunda"
$ws.Cells.Item(56, 13).Value = 300
$ws.Cells.Item(56, 14).Value = 20000
$ws.Cells.Item(56, 15).Value = 21000
$ws.Cells.Item(56, 16).Value = 20500
$ws.Cells.Item(56, 17).Value = "$/caja 18 kilos empedrada"
$ws.Cells.Item(56, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(56, 19).Value = 1139

# Row 57
$ws.Cells.Item(57, 4).Value = 44278
$ws.Cells.Item(57, 11).Value = "Granny Smith"
$ws.Cells.Item(57, 12).Value = "Segunda"
$ws.Cells.Item(57, 13).Value = 250
$ws.Cells.Item(57, 14).Value = 20000
$ws.Cells.Item(57, 15).Value = 21000
$ws.Cells.Item(57, 16).Value = 20500
$ws.Cells.Item(57, 17).Value = "$/caja 18 kilos empedrada"
$ws.Cells.Item(57, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(57, 19).Value = 1139

# Row 58
$ws.Cells.Item(58, 4).Value = 44278
$ws.Cells.Item(58, 11).Value = "Royal Gala"
$ws.Cells.Item(58, 12).Value = "Segunda"
$ws.Cells.Item(58, 13).Value = 270
$ws.Cells.Item(58, 14).Value = 21000
$ws.Cells.Item(58, 15).Value = 22000
$ws.Cells.Item(58, 16).Value = 21500
$ws.Cells.Item(58, 17).Value = "$/caja 18 kilos empedrada"
$ws.Cells.Item(58, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(58, 19).Value = 1194

# Row 59
$ws.Cells.Item(59, 4).Value = 44202
$ws.Cells.Item(59, 11).Value = "Granny Smith"
$ws.Cells.Item(59, 12).Value = "Primera"
$ws.Cells.Item(59, 13).Value = 250
$ws.Cells.Item(59, 14).Value = 25000
$ws.Cells.Item(59, 15).Value = 26000
$ws.Cells.Item(59, 16).Value = 25500
$ws.Cells.Item(59, 17).Value = "$/caja 18 kilos empedrada"
$ws.Cells.Item(59, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(59, 19).Value = 1417

# Row 60
$ws.Cells.Item(60, 4).Value = 44307
$ws.Cells.Item(60, 11).Value = "Fuji royal"
$ws.Cells.Item(60, 12).Value = "Calibre 80"
$ws.Cells.Item(60, 13).Value = 250
$ws.Cells.Item(60, 14).Value = 19000
$ws.Cells.Item(60, 15).Value = 20000
$ws.Cells.Item(60, 16).Value = 19500
$ws.Cells.Item(60, 17).Value = "$/caja 18 kilos embalada"
$ws.Cells.Item(60, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(60, 19).Value = 1083

# Row 61
$ws.Cells.Item(61, 4).Value = 44307
$ws.Cells.Item(61, 11).Value = "Granny Smith"
$ws.Cells.Item(61, 12).Value = "Calibre 80"
$ws.Cells.Item(61, 13).Value = 300
$ws.Cells.Item(61, 14).Value = 19000
$ws.Cells.Item(61, 15).Value = 20000
$ws.Cells.Item(61, 16).Value = 19500
$ws.Cells.Item(61, 17).Value = "$/caja 18 kilos embalada"
$ws.Cells.Item(61, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(61, 19).Value = 1083

# Row 62
$ws.Cells.Item(62, 4).Value = 44307
$ws.Cells.Item(62, 11).Value = "Royal Gala"
$ws.Cells.Item(62, 12).Value = "Calibre 90"
$ws.Cells.Item(62, 13).Value = 250
$ws.Cells.Item(62, 14).Value = 19000
$ws.Cells.Item(62, 15).Value = 20000
$ws.Cells.Item(62, 16).Value = 19500
$ws.Cells.Item(62, 17).Value = "$/caja 18 kilos embalada"
$ws.Cells.Item(62, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(62, 19).Value = 1083

# Row 63
$ws.Cells.Item(63, 4).Value = 44161
$ws.Cells.Item(63, 11).Value = "Fuji royal"
$ws.Cells.Item(63, 12).Value = "Calibre 90"
$ws.Cells.Item(63, 13).Value = 300
$ws.Cells.Item(63, 14).Value = 24000
$ws.Cells.Item(63, 15).Value = 25000
$ws.Cells.Item(63, 16).Value = 24500
$ws.Cells.Item(63, 17).Value = "$/caja 18 kilos embalada"
$ws.Cells.Item(63, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(63, 19).Value = 1361

# Row 64
$ws.Cells.Item(64, 1).Value = $mercadoId
$ws.Cells.Item(64, 2).Value = $mercado
$ws.Cells.Item(64, 3).Value = $region
$ws.Cells.Item(64, 4).Value = 44161
$ws.Cells.Item(64, 4).NumberFormat = $fechaFormat
$ws.Cells.Item(64, 5).Value = $codreg
$ws.Cells.Item(64, 6).Value = $tipo
$ws.Cells.Item(64, 7).Value = $productoId
$ws.Cells.Item(64, 8).Value = $producto
$ws.Cells.Item(64, 9).Value = $categoriaId
$ws.Cells.Item(64, 10).Value = $categoria
$ws.Cells.Item(64, 11).Value = "Granny Smith"
$ws.Cells.Item(64, 12).Value = "Calibre 90"
$ws.Cells.Item(64, 13).Value = 320
$ws.Cells.Item(64, 14).Value = 24000
$ws.Cells.Item(64, 15).Value = 25000
$ws.Cells.Item(64, 16).Value = 24500
$ws.Cells.Item(64, 17).Value = "$/caja 18 kilos embalada"
$ws.Cells.Item(64, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(64, 19).Value = 1361
$ws.Cells.Item(64, 20).Value = $kgUnidad

# Row 65
$ws.Cells.Item(65, 1).Value = $mercadoId
$ws.Cells.Item(65, 2).Value = $mercado
$ws.Cells.Item(65, 3).Value = $region
$ws.Cells.Item(65, 4).Value = 44334
$ws.Cells.Item(65, 4).NumberFormat = $fechaFormat
$ws.Cells.Item(65, 5).Value = $codreg
$ws.Cells.Item(65, 6).Value = $tipo
$ws.Cells.Item(65, 7).Value = $productoId
$ws.Cells.Item(65, 8).Value = $producto
$ws.Cells.Item(65, 9).Value = $categoriaId
$ws.Cells.Item(65, 10).Value = $categoria
$ws.Cells.Item(65, 11).Value = "Fuji royal"
$ws.Cells.Item(65, 12).Value = "Calibre 80"
$ws.Cells.Item(65, 13).Value = 270
$ws.Cells.Item(65, 14).Value = 16000
$ws.Cells.Item(65, 15).Value = 17000
$ws.Cells.Item(65, 16).Value = 16500
$ws.Cells.Item(65, 17).Value = "$/caja 18 kilos embalada"
$ws.Cells.Item(65, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(65, 19).Value = 917
$ws.Cells.Item(65, 20).Value = $kgUnidad

# Row 66
$ws.Cells.Item(66, 1).Value = $mercadoId
$ws.Cells.Item(66, 2).Value = $mercado
$ws.Cells.Item(66, 3).Value = $region
$ws.Cells.Item(66, 4).Value = 44334
$ws.Cells.Item(66, 4).NumberFormat = $fechaFormat
$ws.Cells.Item(66, 5).Value = $codreg
$ws.Cells.Item(66, 6).Value = $tipo
$ws.Cells.Item(66, 7).Value = $productoId
$ws.Cells.Item(66, 8).Value = $producto
$ws.Cells.Item(66, 9).Value = $categoriaId
$ws.Cells.Item(66, 10).Value = $categoria
$ws.Cells.Item(66, 11).Value = "Granny Smith"
$ws.Cells.Item(66, 12).Value = "Calibre 90"
$ws.Cells.Item(66, 13).Value = 300
$ws.Cells.Item(66, 14).Value = 16000
$ws.Cells.Item(66, 15).Value = 17000
$ws.Cells.Item(66, 16).Value = 16500
$ws.Cells.Item(66, 17).Value = "$/caja 18 kilos embalada"
$ws.Cells.Item(66, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(66, 19).Value = 917
$ws.Cells.Item(66, 20).Value = $kgUnidad

# Row 67
$ws.Cells.Item(67, 1).Value = $mercadoId
$ws.Cells.Item(67, 2).Value = $mercado
$ws.Cells.Item(67, 3).Value = $region
$ws.Cells.Item(67, 4).Value = 44334
$ws.Cells.Item(67, 4).NumberFormat = $fechaFormat
$ws.Cells.Item(67, 5).Value = $codreg
$ws.Cells.Item(67, 6).Value = $tipo
$ws.Cells.Item(67, 7).Value = $productoId
$ws.Cells.Item(67, 8).Value = $producto
$ws.Cells.Item(67, 9).Value = $categoriaId
$ws.Cells.Item(67, 10).Value = $categoria
$ws.Cells.Item(67, 11).Value = "Royal Gala"
$ws.Cells.Item(67, 12).Value = "Calibre 90"
$ws.Cells.Item(67, 13).Value = 300
$ws.Cells.Item(67, 14).Value = 16000
$ws.Cells.Item(67, 15).Value = 17000
$ws.Cells.Item(67, 16).Value = 16500
$ws.Cells.Item(67, 17).Value = "$/caja 18 kilos embalada"
$ws.Cells.Item(67, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(67, 19).Value = 917
$ws.Cells.Item(67, 20).Value = $kgUnidad

# Row 68
$ws.Cells.Item(68, 1).Value = $mercadoId
$ws.Cells.Item(68, 2).Value = $mercado
$ws.Cells.Item(68, 3).Value = $region
$ws.Cells.Item(68, 4).Value = 44432
$ws.Cells.Item(68, 4).NumberFormat = $fechaFormat
$ws.Cells.Item(68, 5).Value = $codreg
$ws.Cells.Item(68, 6).Value = $tipo
$ws.Cells.Item(68, 7).Value = $productoId
$ws.Cells.Item(68, 8).Value = $producto
$ws.Cells.Item(68, 9).Value = $categoriaId
$ws.Cells.Item(68, 10).Value = $categoria
$ws.Cells.Item(68, 11).Value = "Granny Smith"
$ws.Cells.Item(68, 12).Value = "Calibre 90"
$ws.Cells.Item(68, 13).Value = 300
$ws.Cells.Item(68, 14).Value = 17000
$ws.Cells.Item(68, 15).Value = 18000
$ws.Cells.Item(68, 16).Value = 17500
$ws.Cells.Item(68, 17).Value = "$/caja 18 kilos empedrada"
$ws.Cells.Item(68, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(68, 19).Value = 972
$ws.Cells.Item(68, 20).Value = $kgUnidad

# Row 69
$ws.Cells.Item(69, 1).Value = $mercadoId
$ws.Cells.Item(69, 2).Value = $mercado
$ws.Cells.Item(69, 3).Value = $region
$ws.Cells.Item(69, 4).Value = 44432
$ws.Cells.Item(69, 4).NumberFormat = $fechaFormat
$ws.Cells.Item(69, 5).Value = $codreg
$ws.Cells.Item(69, 6).Value = $tipo
$ws.Cells.Item(69, 7).Value = $productoId
$ws.Cells.Item(69, 8).Value = $producto
$ws.Cells.Item(69, 9).Value = $categoriaId
$ws.Cells.Item(69, 10).Value = $categoria
$ws.Cells.Item(69, 11).Value = "Pink Lady"
$ws.Cells.Item(69, 12).Value = "Calibre 80"
$ws.Cells.Item(69, 13).Value = 250
$ws.Cells.Item(69, 14).Value = 17000
$ws.Cells.Item(69, 15).Value = 18000
$ws.Cells.Item(69, 16).Value = 17500
$ws.Cells.Item(69, 17).Value = "$/caja 18 kilos empedrada"
$ws.Cells.Item(69, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(69, 19).Value = 972
$ws.Cells.Item(69, 20).Value = $kgUnidad

# Row 70
$ws.Cells.Item(70, 1).Value = $mercadoId
$ws.Cells.Item(70, 2).Value = $mercado
$ws.Cells.Item(70, 3).Value = $region
$ws.Cells.Item(70, 4).Value = 44432
$ws.Cells.Item(70, 4).NumberFormat = $fechaFormat
$ws.Cells.Item(70, 5).Value = $codreg
$ws.Cells.Item(70, 6).Value = $tipo
$ws.Cells.Item(70, 7).Value = $productoId
$ws.Cells.Item(70, 8).Value = $producto
$ws.Cells.Item(70, 9).Value = $categoriaId
$ws.Cells.Item(70, 10).Value = $categoria
$ws.Cells.Item(70, 11).Value = "Royal Gala"
$ws.Cells.Item(70, 12).Value = "Calibre 90"
$ws.Cells.Item(70, 13).Value = 520
$ws.Cells.Item(70, 14).Value = 17000
$ws.Cells.Item(70, 15).Value = 18000
$ws.Cells.Item(70, 16).Value = 17500
$ws.Cells.Item(70, 17).Value = "$/caja 18 kilos empedrada"
$ws.Cells.Item(70, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(70, 19).Value = 972
$ws.Cells.Item(70, 20).Value = $kgUnidad

# Row 71
$ws.Cells.Item(71, 1).Value = $mercadoId
$ws.Cells.Item(71, 2).Value = $mercado
$ws.Cells.Item(71, 3).Value = $region
$ws.Cells.Item(71, 4).Value = 44432
$ws.Cells.Item(71, 4).NumberFormat = $fechaFormat
$ws.Cells.Item(71, 5).Value = $codreg
$ws.Cells.Item(71, 6).Value = $tipo
$ws.Cells.Item(71, 7).Value = $productoId
$ws.Cells.Item(71, 8).Value = $producto
$ws.Cells.Item(71, 9).Value = $categoriaId
$ws.Cells.Item(71, 10).Value = $categoria
$ws.Cells.Item(71, 11).Value = "Scarlett"
$ws.Cells.Item(71, 12).Value = "Calibre 80"
$ws.Cells.Item(71, 13).Value = 300
$ws.Cells.Item(71, 14).Value = 17000
$ws.Cells.Item(71, 15).Value = 18000
$ws.Cells.Item(71, 16).Value = 17500
$ws.Cells.Item(71, 17).Value = "$/caja 18 kilos empedrada"
$ws.Cells.Item(71, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(71, 19).Value = 972
$ws.Cells.Item(71, 20).Value = $kgUnidad
